$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Image, title, and sentence: express" + " how I'm unique"
#    -> merge into a single run's text (About page bullet).
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("express how I’m unique", $true, $false, $false, $false, $false, $true, 1, $false, "express how I’m unique", 2)

# ---------------------------------------------------------------------
# 2) "Image, title, and sentence: " + "software engineer, what kind of developer I am"
#    -> merge into a single run's text (Portfolio page bullet).
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("Image, title, and sentence: software engineer", $true, $false, $false, $false, $false, $true, 1, $false, "Image, title, and sentence: software engineer", 2)

# ---------------------------------------------------------------------
# 3) "Github" + " account" (wrapped in proofErr spellStart/spellEnd)
#    -> "Git" / "H" / "ub account" as three clean runs, no proofErr.
# ---------------------------------------------------------------------
$rngFind = $d.Content
$gotIt = $rngFind.Find.Execute("Github account")
$startPos = $rngFind.Start

# Insert a brand-new (proofErr-free) paragraph right before the old one;
# Word clones the paragraph formatting (style/numbering/rPr) automatically.
$ins = $d.Range($startPos, $startPos)
$ins.InsertParagraphBefore()

# Type the replacement text into the new paragraph.
$ins2 = $d.Range($startPos, $startPos)
$ins2.InsertAfter("GitHub account")

# Force a run split after "Git" and after "GitH" by toggling a trivial
# character-format property on a 1-character range and back again; Word
# splits runs at format boundaries, and reverting the property keeps the
# visible formatting identical across the new runs.
$b1 = $d.Range($startPos + 3, $startPos + 4)
$b1.Font.Bold = 1
$b1.Font.Bold = 0

# Remove the old paragraph (its text plus the paragraph mark), taking the
# stale proofErr markers with it.
$oldParaStart = $startPos + 15
$oldPara = $d.Range($oldParaStart, $oldParaStart + 15)
$oldPara.Delete()

# ---------------------------------------------------------------------
# 4) "Color:" gains a trailing run " light gray, dark blue"
# ---------------------------------------------------------------------
$rngColor = $d.Content
$rngColor.Find.Execute("Color:") | Out-Null
$rngColor.Collapse(0)
$rngColor.InsertAfter(" light gray, dark blue")
